$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: existing transaction updated (farmer re-keyed, amount corrected) ---
$ws.Range("A2").Value = "RameshPatil"
$ws.Range("B2").Value = "'9090909090"
$ws.Range("B2").ClearFormats()
$ws.Range("D2").Value = 280

# --- Row 3 (new): same farmer, earlier transaction ---
$ws.Range("A3").Value = "RameshPatil"
$ws.Range("B3").Value = "'9090909090"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = "'2025-03-17"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = 278
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "Cow"

# --- Row 4 (new): another farmer's transaction ---
$ws.Range("A4").Value = "Rohit"
$ws.Range("B4").Value = "'9191919191"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "'2025-04-17"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = 369.6
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = "Cow"
